# "Generate Report for Handoff" — refresh the localization-status report.
#
# The handback/handoff polling job re-ran and produced new, unified
# timestamps for the files that previously showed stale/mismatched
# datetimes across the Overview, zh-cn and de-de sheets. This mirrors
# that refresh by writing the newly observed timestamp into every cell
# that held one of the old (now-superseded) timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ---
$overview = $wb.Worksheets.Item("Overview")
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $overview.Cells.Item($r, 4).Value = "2016-03-21 20:28:10"
}

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $zhcnRows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-21 20:28:06"
}

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$dede = $wb.Worksheets.Item("de-de")
$dedeRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $dedeRows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-21 20:28:10"
}
